$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the "through" date references
$ws.Name = "Through 2021-12-31"
$ws.Range("A13").Value = "December (through 12-31)"

# Update H12
$ws.Range("H12").Value = 201

# Update row 13 (December data)
$ws.Range("B13").Value = 49
$ws.Range("C13").Value = 100
$ws.Range("D13").Value = 116
$ws.Range("E13").Value = 82
$ws.Range("F13").Value = 69
$ws.Range("G13").Value = 149
$ws.Range("H13").Value = 204

# Update row 14 (Total)
$ws.Range("B14").Value = 340
$ws.Range("C14").Value = 663
$ws.Range("D14").Value = 937
$ws.Range("E14").Value = 764
$ws.Range("F14").Value = 603
$ws.Range("G14").Value = 1413
$ws.Range("H14").Value = 1848
